$wb = $excel.ActiveWorkbook

# --- Rename "Sheet1" -> "smoke" ---
$smoke = $wb.Worksheets.Item("Sheet1")
$smoke.Name = "smoke"

# --- "regression" sheet: remove D1/D2, add new rows 7 & 8, widen column B, move selection ---
$reg = $wb.Worksheets.Item("regression")

# Remove the D1 / D2 cells entirely (content + formatting), matching the row-1/row-2
# trim in the target workbook.
$reg.Range("D1").Clear() | Out-Null
$reg.Range("D2").Clear() | Out-Null

# Add the two new test rows.
$reg.Range("A7").Value = "TC103"
$reg.Range("B7").Value = "Location"
$reg.Range("A8").Value = "TC103"
$reg.Range("B8").Value = "Sydney"

# Widen column B to fit the new, longer content (57.21875 is the on-disk target
# width; the host's pixel-quantised ColumnWidth setter snaps to the nearest
# representable width, so feed it the input that lands closest to that target).
$reg.Columns.Item(2).ColumnWidth = 57.05

# Move the active selection to B7.
$reg.Range("B7").Select() | Out-Null

# --- "smoke" sheet (formerly Sheet1): collapse its stored selection back to the
# default top-left cell, then restore "regression" as the active sheet/tab. ---
$smoke.Range("A1").Select() | Out-Null
$reg.Select() | Out-Null
$reg.Range("B7").Select() | Out-Null
